$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Popisek" (description) texts typed into column G (Teorie block) ---
# Typed in roughly this order so the rebuilt shared-string table matches the
# order new unique strings were introduced in the edited workbook.
$ws.Range("G3").Value = "tohle je lehčí než zobrazení"
$ws.Range("G8").Value = "praktické téma, docela jde"
$ws.Range("G9").Value = "podobné 1 ale je tam par věci  co nejsou zastak jednoduchy"
$ws.Range("G10").Value = "proste tohle nam říká zda mame vůbec reseni"
$ws.Range("G11").Value = "je to nic moc ale není to až tak dlouhy"
$ws.Range("G12").Value = "docela jde, něco vim prakticky"
$ws.Range("G13").Value = "podobné 11 ale mozna trosku lehci"
$ws.Range("G14").Value = "jednoduší než téma 7"
$ws.Range("G15").Value = "není to ažtak obtížné"
$ws.Range("G16").Value = "tohle jde"
$ws.Range("G17").Value = "15 tema bylo lepsi"
$ws.Range("G18").Value = "narocnejsi kvuli množství textu "
$ws.Range("G19").Value = "asi nejlehčí téma"
$ws.Range("G20").Value = "docela kratky, jde to"
$ws.Range("G21").Value = "ty priklady hodne pomahaji"
$ws.Range("C7").Value = "trochu se to podobá 1"
$ws.Range("C8").Value = "nedelat chyby v GEM"
$ws.Range("G22").Value = "vzorců tam je extrémně moc"
$ws.Range("G23").Value = "extrém"
$ws.Range("G24").Value = "prakticky trochu podobné tématu 11 a navazuje to na tema 22"
$ws.Range("C5").Value = "jde to"
$ws.Range("G26").Value = "druhy nejlehci"
$ws.Range("G2").Value = "neumim to rict presne podle tech vet nektery věci"
$ws.Range("G6").Value = "hodne v poho"
$ws.Range("G7").Value = "jde to, jsou tam věci co se daji snadno zapomenout"
$ws.Range("G5").Value = "jde to, je to v poho"
$ws.Range("G4").Value = "tohle je dobry, krome toho co znamena Podprostor"

# --- Cells that reuse an already-existing text (no new unique string) ---
$ws.Range("C1").Value = "Obtížnost/Popisek"
$ws.Range("G1").Value = "Obtížnost/Popisek"
$ws.Range("C6").Value = "dá se v tom chybovat"
$ws.Range("G25").Value = "jde to"

# --- Updated "Znamka" (grade) numbers ---
$ws.Range("B5").Value = 1.2
$ws.Range("B6").Value = 3.4
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3.4

$ws.Range("F2").Value = 3.4
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2.3
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 1.2
$ws.Range("F7").Value = 2.3
$ws.Range("F8").Value = 1.2
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 3.4
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 3.4
$ws.Range("F13").Value = 3.4
$ws.Range("F14").Value = 2.3
$ws.Range("F15").Value = 2.3
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 3.4
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 1.2
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 2
$ws.Range("F22").Value = 3.4
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 1.2

# --- Column G got much wider (long free-text descriptions); H/I slightly narrower ---
$ws.Columns.Item(7).ColumnWidth = 175
$ws.Range("H:I").ColumnWidth = 8.6

# --- View state: zoomed out a bit, selection moved to G8 ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("G8").Select()
